$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 139, pushing the existing rows 139..232 down to 140..233.
$ws.Rows(139).Insert()

# Seed the new row with a copy of the (now shifted) row 140 so every column
# (including formatting) matches the surrounding data rows.
$ws.Range("A140:T140").Copy()
$ws.Range("A139").PasteSpecial()

# Overwrite the columns that differ for this new record:
# Fecha (D), Precio mínimo (N), Precio máximo (O), Precio promedio (P), Precio promedio $/kilo (S)
$ws.Cells.Item(139, 4).Value = 45126
$ws.Cells.Item(139, 14).Value = 9000
$ws.Cells.Item(139, 15).Value = 10000
$ws.Cells.Item(139, 16).Value = 9500
$ws.Cells.Item(139, 19).Value = 950
